# feat: update Time Reward config / 更新在线奖励表
#
# Rebalance the "time needed" (column C, seconds) and the reward payload
# string (column E, "gold|diamond|petId|activityToken") for every online
# reward tier (rows 5-16) on the only worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C: seconds of online time required for each gift tier.
$ws.Range("C5").Value = 300
$ws.Range("C6").Value = 600
$ws.Range("C8").Value = 1800
$ws.Range("C9").Value = 3600
$ws.Range("C10").Value = 7200
$ws.Range("C11").Value = 10800
$ws.Range("C12").Value = 18000
$ws.Range("C13").Value = 28800
$ws.Range("C14").Value = 43200
$ws.Range("C15").Value = 61200
$ws.Range("C16").Value = 86400

# Column E: reward payload strings, written in the same order the rows
# were reworked so the new entries land where the refreshed table expects
# them (rows 10-16 first, then 5, 7, 9, 8, 6).
$ws.Range("E10").Value = "15000|1000|0|0"
$ws.Range("E11").Value = "20000|1500|0|0"
$ws.Range("E12").Value = "30000|2000|0|0"
$ws.Range("E13").Value = "40000|2500|0|0"
$ws.Range("E14").Value = "50000|3000|0|0"
$ws.Range("E15").Value = "50000|4000|0|0"
$ws.Range("E16").Value = "50000|5000|0|0"
$ws.Range("E5").Value = "500|100|0|0"
$ws.Range("E7").Value = "3000|300|0|0"
$ws.Range("E9").Value = "10000|600|0|0"
$ws.Range("E8").Value = "6000|500|0|0"
$ws.Range("E6").Value = "1000|200|0|0"

# The author's saved cursor position ends on G8.
$ws.Range("G8").Select()
